$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ParameterA = 0.120"
$ws.Range("A3").Value = "ParameterB = 0.050"
$ws.Range("A4").Value = "ParameterC = -0.100"
$ws.Range("A5").Value = "ParameterD = -0.100"
$ws.Range("A6").Value = "ParameterE = 0.000"
$ws.Range("A7").Value = "ParameterF = 0.000"
$ws.Range("A8").Value = "MaxDD = 0.338"
$ws.Range("A9").Value = "NetProfit = 34504303.6"
$ws.Range("A10").Value = "SharpeRatio = 1.389"
$ws.Range("A11").Value = "AnnualizedReturn = 0.998"
